# Refresh the cryptos list: update Price (D) and Volume(1h) (E) columns.
# Values are written as text (matching the original inline-string cells),
# so numeric-looking prices (e.g. "213.70", "27.332.72") keep their exact
# formatting instead of being reinterpreted as numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.332.72'
$c.Style = $s

$c = $ws.Range('E2')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.14%  '
$c.Style = $s

$c = $ws.Range('D3')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.647.82'
$c.Style = $s

$c = $ws.Range('E3')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.18%  '
$c.Style = $s

$c = $ws.Range('E4')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.Style = $s

$c = $ws.Range('D5')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '213.70'
$c.Style = $s

$c = $ws.Range('E5')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.89%  '
$c.Style = $s

$c = $ws.Range('E6')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.35%  '
$c.Style = $s

$c = $ws.Range('E7')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.Style = $s

$c = $ws.Range('D8')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '23.98'
$c.Style = $s

$c = $ws.Range('E8')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = $s

$c = $ws.Range('E9')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.27%  '
$c.Style = $s

$c = $ws.Range('E10')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.27%  '
$c.Style = $s

$c = $ws.Range('D11')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0877'
$c.Style = $s

$c = $ws.Range('E11')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.36%  '
$c.Style = $s

$c = $ws.Range('D12')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.881.42'
$c.Style = $s

$c = $ws.Range('E12')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.40%  '
$c.Style = $s

$c = $ws.Range('D13')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.649.45'
$c.Style = $s

$c = $ws.Range('E13')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '
$c.Style = $s

$c = $ws.Range('D14')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.08'
$c.Style = $s

$c = $ws.Range('E14')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.48%  '
$c.Style = $s

$c = $ws.Range('D15')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.567'
$c.Style = $s

$c = $ws.Range('E15')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  +1.63%  '
$c.Style = $s

$c = $ws.Range('D16')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '65.73'
$c.Style = $s

$c = $ws.Range('E16')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.14%  '
$c.Style = $s

$c = $ws.Range('D17')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.346.86'
$c.Style = $s

$c = $ws.Range('E17')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.08%  '
$c.Style = $s

$c = $ws.Range('D18')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '234.04'
$c.Style = $s

$c = $ws.Range('E18')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -6.84%  '
$c.Style = $s

$c = $ws.Range('D19')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0₃0724'
$c.Style = $s

$c = $ws.Range('E19')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.63%  '
$c.Style = $s

$c = $ws.Range('D20')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.43'
$c.Style = $s

$c = $ws.Range('E20')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.36%  '
$c.Style = $s

$c = $ws.Range('E21')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = $s

$c = $ws.Range('D22')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.40'
$c.Style = $s

$c = $ws.Range('E22')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.09%  '
$c.Style = $s

$c = $ws.Range('D23')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.25'
$c.Style = $s

$c = $ws.Range('E23')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.19%  '
$c.Style = $s

$c = $ws.Range('E24')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.44%  '
$c.Style = $s

$c = $ws.Range('D25')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '146.14'
$c.Style = $s

$c = $ws.Range('E25')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.98%  '
$c.Style = $s

$c = $ws.Range('E26')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.31%  '
$c.Style = $s

$c = $ws.Range('E27')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.65%  '
$c.Style = $s

$c = $ws.Range('E28')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c.Style = $s

$c = $ws.Range('D29')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.Style = $s

$c = $ws.Range('E29')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.49%  '
$c.Style = $s

$c = $ws.Range('E30')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.36%  '
$c.Style = $s

$c = $ws.Range('E31')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '
$c.Style = $s

$c = $ws.Range('E32')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.48%  '
$c.Style = $s

$c = $ws.Range('D33')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.456.22'
$c.Style = $s

$c = $ws.Range('E33')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.02%  '
$c.Style = $s

$c = $ws.Range('E34')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.15%  '
$c.Style = $s

$c = $ws.Range('E35')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -4.57%  '
$c.Style = $s

$c = $ws.Range('E36')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.Style = $s

$c = $ws.Range('D37')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.907'
$c.Style = $s

$c = $ws.Range('E37')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -5.22%  '
$c.Style = $s

$c = $ws.Range('D38')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.571'
$c.Style = $s

$c = $ws.Range('E38')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.85%  '
$c.Style = $s

$c = $ws.Range('E39')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.28%  '
$c.Style = $s

$c = $ws.Range('E40')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c.Style = $s

$c = $ws.Range('E41')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c.Style = $s

$c = $ws.Range('D42')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '65.20'
$c.Style = $s

$c = $ws.Range('E42')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -5.72%  '
$c.Style = $s

$c = $ws.Range('E43')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.60%  '
$c.Style = $s

$c = $ws.Range('E44')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.61%  '
$c.Style = $s

$c = $ws.Range('D45')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.789.53'
$c.Style = $s

$c = $ws.Range('E45')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -3.41%  '
$c.Style = $s

$c = $ws.Range('D46')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.781'
$c.Style = $s

$c = $ws.Range('E46')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.28%  '
$c.Style = $s

$c = $ws.Range('E47')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.Style = $s

$c = $ws.Range('D48')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '88.13'
$c.Style = $s

$c = $ws.Range('E48')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.Style = $s

$c = $ws.Range('E49')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -0.93%  '
$c.Style = $s

$c = $ws.Range('E50')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -1.97%  '
$c.Style = $s

$c = $ws.Range('D51')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.77'
$c.Style = $s

$c = $ws.Range('E51')
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '  -2.13%  '
$c.Style = $s

